$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "capacidad" column (D) so capacidad moves to E
$ws.Columns.Item(4).Insert()

# Update header row: C1 = disponibilidad (replacing "estado"), E1 = capacidad
$ws.Range("C1").Value = "disponibilidad"
$ws.Range("E1").Value = "capacidad"

# Underline a stray far-away cell first (so its style claims the next style slot
# before the temporary quote-prefixed style used below is created)
$ws.Range("D10").Font.Underline = 1

# Fill new "activo" column with literal text "TRUE" for each data row
$ws.Range("D2").Value = "'TRUE"
$ws.Range("D3").Value = "'TRUE"
$ws.Range("D4").Value = "'TRUE"
$ws.Range("D5").Value = "'TRUE"

# Normalize formatting on the new cells to match the rest of the table
# (remove the quote-prefix styling introduced by the text entry above)
$ws.Range("C2").Copy()
$ws.Range("D2:D5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Now set the activo header (added after the TRUE values, matching shared string order)
$ws.Range("D1").Value = "activo"

# Ensure page setup reflects default portrait orientation
$ws.PageSetup.Orientation = 1

# Leave the selection further down-right, matching the final interactive state of the sheet
[void]$ws.Range("G11").Select()

$wb.Save()
